# "only show not closed project"
#
# - Localization!C3  : Project "HVS 1BC" -> "HVS 1BC 2"   (project renamed / re-flagged)
# - Energy!C3        : Project "HVS 1BC" -> "HVS 1BC 2"
# - Energy!B5        : SIE "Sateria Salim" -> "someone"
# - Energy!M5        : Biz award date 7/31/2024 -> 7/19/2024
# - Energy!W5        : Action item for Cindy updated text
# - Localization tab becomes the active/selected tab (was Energy)

$wb = $excel.ActiveWorkbook

$wsLocalization = $wb.Worksheets.Item("Localization")
$wsOthers       = $wb.Worksheets.Item("Others")
$wsEnergy       = $wb.Worksheets.Item("Energy")

# --- Data edits ------------------------------------------------------------
# Order matches the shared-string table append order produced by the
# original edit (someone -> action itemenergy... -> HVS 1BC 2).

$wsEnergy.Range("B5").Value = "someone"
$wsEnergy.Range("M5").Value = "7/19/2024"
$wsEnergy.Range("W5").Value = "action itemenergy aasdfas for cindy 4"

$wsLocalization.Range("C3").Value = "HVS 1BC 2"
$wsEnergy.Range("C3").Value = "HVS 1BC 2"

# --- View / selection state --------------------------------------------

# Energy was the active tab with C3 selected; move the selection there first
# so it's no longer the last-activated sheet once Localization is activated
# below.
$wsEnergy.Activate()
$wsEnergy.Range("T5:X7").Select()

# Localization becomes the active (selected) sheet / tab, with C3 selected.
$wsLocalization.Activate()
$wsLocalization.Range("C3").Select()
